$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Significant Components": reorder the variable names listed inside
# the significant-components text cells (C2, C3, C5, C6) to match the new
# factor-variable ordering produced upstream.
# ---------------------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significant Components")
$wsSig.Range("C2").Value = "['QESL' 'PPUNIT' 'QEDLESHI' 'QHISPC' 'QNOHLTH' 'QSERV' 'QEXTRCT' 'PERCAP'`n 'QFHH']"
$wsSig.Range("C3").Value = "['PERCAP' 'QRICH' 'MDHSEVAL']"
$wsSig.Range("C5").Value = "['QRENTER' 'QNOAUTO' 'QPOVTY']"
$wsSig.Range("C6").Value = "['QAGEDEP' 'QFEMALE' 'QFEMLBR']"

# ---------------------------------------------------------------------------
# Sheet "Loading Factors": the variable rows (A2:F20) are re-sorted into the
# new order, and the loading values are refreshed with the (very slightly
# different, re-computed) figures that came with the reorder.
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Loading Factors")

$loadRows = @(
    @(2, "QESL", 0.7710373703114096, 0.1699122220342584, -0.04922328053492717, 0.1763260275555902, -0.2279019739325974),
    @(3, "PPUNIT", 0.7163381413134868, -0.04911535347006139, -0.08826232351490147, -0.3752405587914745, 0.1082069730305032),
    @(4, "QEDLESHI", 0.8620380840342905, 0.2178973314144597, 0.00913724299824761, 0.1896693359116807, -0.1047063828551197),
    @(5, "QHISPC", 0.8195686393389251, 0.3528442964268741, -0.1123103918846655, 0.1083518371166312, -0.1312378922779179),
    @(6, "QNOHLTH", 0.664640440656863, 0.4300115824886732, -0.06979631800196752, 0.2528269903147914, -0.1342238767271383),
    @(7, "QSERV", 0.5739901519602962, 0.3660124165492233, -0.1657718182772126, 0.3037369625029354, -0.05431698638550993),
    @(8, "QEXTRCT", 0.7782826624132134, 0.1356252026708027, -0.02831620240160483, 0.06461226540028045, -0.2142235158524506),
    @(9, "PERCAP", 0.5008581706111616, 0.706828701298429, -0.2354656060493542, 0.1849142842472271, 0.08409873372830448),
    @(10, "QFHH", 0.5756381644391316, 0.2481999142322477, -0.007302870366954998, 0.08150148524366285, 0.2295452304581432),
    @(11, "QRICH", 0.2370191025212669, 0.837586213503246, -0.2059104174342114, 0.3148028664732525, -0.03073101713398335),
    @(12, "MDHSEVAL", 0.3768770853591332, 0.7947218037000914, -0.07404588234023969, -0.03403758781339291, 0.0291764094771565),
    @(13, "MEDAGE", -0.2799970105938261, -0.2283516395058378, 0.7796766778756494, -0.3080563892774956, -0.07828734457891075),
    @(14, "QAGEDEP", -0.01911179177596987, -0.1421335481267801, 0.722804315970739, -0.06893936282044803, 0.5773032568375647),
    @(15, "QSSBEN", 0.03069875454850582, -0.06832501990930143, 0.8135614489793346, -0.1401439056242418, 0.1110613116324732),
    @(16, "QRENTER", -0.02252423251283685, 0.2520805273154729, -0.4372443891608503, 0.737637042952187, -0.1126483035129214),
    @(17, "QNOAUTO", 0.133889969887876, 0.0532923984052213, -0.04656029707767763, 0.7398553445767126, 0.01303061619987351),
    @(18, "QPOVTY", 0.2772247094325571, 0.1320875694826831, -0.330793361418331, 0.5596088471642952, 0.1365168950783244),
    @(19, "QFEMALE", -0.02635136390715215, -0.04008883933396452, 0.2036171317561216, 0.0233369804047918, 0.855443602512891),
    @(20, "QFEMLBR", -0.2223644614709627, 0.08523116677994901, -0.0413046652323644, 0.004904383437521533, 0.7515631002766183)
)

foreach ($row in $loadRows) {
    $r = $row[0]
    $wsLoad.Cells.Item($r, 1).Value = $row[1]
    $wsLoad.Cells.Item($r, 2).Value = $row[2]
    $wsLoad.Cells.Item($r, 3).Value = $row[3]
    $wsLoad.Cells.Item($r, 4).Value = $row[4]
    $wsLoad.Cells.Item($r, 5).Value = $row[5]
    $wsLoad.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "All Refactor Variances": the F2 (rightmost 5) columns get refreshed
# with the re-computed variance figures (columns N:R).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Refactor Variances")

$allVarianceUpdates = @(
    @("N2", 4.90494323317218),
    @("O2", 2.604591482109694),
    @("P2", 2.297607997005125),
    @("Q2", 2.045575240043453),
    @("R2", 1.900121960958007),
    @("N3", 0.2581549070090621),
    @("O3", 0.1370837622162997),
    @("P3", 0.1209267366844803),
    @("Q3", 0.1076618547391291),
    @("R3", 0.1000064189977899),
    @("N4", 0.2581549070090621),
    @("O4", 0.3952386692253618),
    @("P4", 0.5161654059098421),
    @("Q4", 0.6238272606489712),
    @("R4", 0.723833679646761),
    @("N5", 0.3566494821504362),
    @("O5", 0.1893857194973272),
    @("P5", 0.1670642581089814),
    @("Q5", 0.1487383880668129),
    @("R5", 0.1381621521764421)
)

foreach ($u in $allVarianceUpdates) {
    $wsAll.Range($u[0]).Value = $u[1]
}

# ---------------------------------------------------------------------------
# Sheet "Final Variances": mirrors the last 5 (F2) columns of
# "All Refactor Variances" above, so it gets the same refreshed figures.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Final Variances")

$finalVarianceUpdates = @(
    @("B2", 4.90494323317218),
    @("C2", 2.604591482109694),
    @("D2", 2.297607997005125),
    @("E2", 2.045575240043453),
    @("F2", 1.900121960958007),
    @("B3", 0.2581549070090621),
    @("C3", 0.1370837622162997),
    @("D3", 0.1209267366844803),
    @("E3", 0.1076618547391291),
    @("F3", 0.1000064189977899),
    @("B4", 0.2581549070090621),
    @("C4", 0.3952386692253618),
    @("D4", 0.5161654059098421),
    @("E4", 0.6238272606489712),
    @("F4", 0.723833679646761),
    @("B5", 0.3566494821504362),
    @("C5", 0.1893857194973272),
    @("D5", 0.1670642581089814),
    @("E5", 0.1487383880668129),
    @("F5", 0.1381621521764421)
)

foreach ($u in $finalVarianceUpdates) {
    $wsFinal.Range($u[0]).Value = $u[1]
}

# ---------------------------------------------------------------------------
# Sheet "Included and Excluded": reorder the variable list in B2 to match
# the new factor-variable ordering (the exclude list in C2 is unchanged).
# ---------------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Included and Excluded")
$wsInc.Range("B2").Value = "[['QESL', 'PPUNIT', 'QEDLESHI', 'QHISPC', 'QNOHLTH', 'QSERV', 'QEXTRCT', 'PERCAP', 'QFHH', 'QRICH', 'MDHSEVAL', 'MEDAGE', 'QAGEDEP', 'QSSBEN', 'QRENTER', 'QNOAUTO', 'QPOVTY', 'QFEMALE', 'QFEMLBR']]"
